# Edit: new linear equations for care-seeking
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: prob_early_ptb -> update B8, add C8 ---
$ws.Range("B8").Value = 0.01
$ws.Range("C8").Value = 0.023

# --- Row 11: prob_late_ptb -> update B11, add C11 ---
$ws.Range("B11").Value = 0.05
$ws.Range("C11").Value = 0.1

# --- Capture the existing comment on D51 before shifting rows ---
$oldComment = $ws.Range("D51").Comment
$commentText = $oldComment.Text()
$oldComment.Delete()

# --- Insert 16 new rows before row 51 (rows 51..66 become new) ---
$ws.Range("A51:A66").EntireRow.Insert()

# --- Populate the 16 new rows with the new care-seeking parameters ---
$ws.Range("A51").Value = "odds_deliver_in_health_centre"
$ws.Range("B51").Value = 0.67

$ws.Range("A52").Value = "rrr_hc_delivery_age_25_29"
$ws.Range("B52").Value = 0.59

$ws.Range("A53").Value = "rrr_hc_delivery_age_30_34"
$ws.Range("B53").Value = 0.27

$ws.Range("A54").Value = "rrr_hc_delivery_age_35_39"
$ws.Range("B54").Value = 0.13

$ws.Range("A55").Value = "rrr_hc_delivery_age_40_44"
$ws.Range("B55").Value = 0.06

$ws.Range("A56").Value = "rrr_hc_delivery_age_45_49"
$ws.Range("B56").Value = 0.02

$ws.Range("A57").Value = "rrr_hc_delivery_married"
$ws.Range("B57").Value = 2.36

$ws.Range("A58").Value = "rrr_hc_delivery_parity_3_to_4"
$ws.Range("B58").Value = 2.21

$ws.Range("A59").Value = "rrr_hc_delivery_parity_>4"
$ws.Range("B59").Value = 4.36

$ws.Range("A60").Value = "rrr_hc_delivery_rural"
$ws.Range("B60").Value = 1.48

$ws.Range("A61").Value = "odds_deliver_at_home"
$ws.Range("B61").Value = 0.06

$ws.Range("A62").Value = "rrr_hb_delivery_age_35_39"
$ws.Range("B62").Value = 0.29

$ws.Range("A63").Value = "rrr_hb_delivery_age_40_44"
$ws.Range("B63").Value = 0.15

$ws.Range("A64").Value = "rrr_hb_delivery_age_45_49"
$ws.Range("B64").Value = 0.04

$ws.Range("A65").Value = "rrr_hb_delivery_parity_3_to_4"
$ws.Range("B65").Value = 3.3

$ws.Range("A66").Value = "rrr_hb_delivery_parity_>4"
$ws.Range("B66").Value = 9

# --- Re-add the comment at its new location, D67 ---
$newComment = $ws.Range("D67").AddComment($commentText)

# --- Fix row 80 (old row 64 "dummy_prob_health_centre") value change 0.7 -> 0.4 ---
$ws.Range("B80").Value = 0.4

# --- Update sheet view to match new authoring position ---
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("A69").Select()

# --- Update calcId to mark workbook as recalculated by a real Excel session ---
$wb.CalcId = 162913
